# Updated symbol list with refreshed Price (D) and Volume(1h) (E) values
# for the cryptos sheet, matching the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" column (D) values, keyed by row number.
$dData = @{
    2 = "299.15"
    3 = "31.65"
    4 = "5.142"
    5 = "0.08075"
    6 = "2.506"
    7 = "7.804"
    8 = "3.905"
    9 = "0.9306"
    10 = "0.1758"
    11 = "0.07415"
    12 = "0.08858"
    13 = "0.03000"
    14 = "0.1000"
    16 = "0.005951"
    17 = "3.528"
    18 = "2.289"
    20 = "0.1338"
    21 = "4.171"
    22 = "0.1680"
    23 = "0.04621"
    25 = "0.004528"
    27 = "0.0003408"
    39 = "0.01753"
    40 = "0.04593"
    41 = "0.006918"
    42 = "0.1373"
    43 = "0.002191"
    44 = "0.01030"
    45 = "0.00006125"
    46 = "0.00000000750"
    47 = "0.008395"
    49 = "0.00002101"
    50 = "0.0002001"
}

# New "Volume(1h)" column (E) values, keyed by row number.
$eData = @{
    2 = "-0.79%"
    3 = "1.00%"
    4 = "0.05%"
    5 = "9.21%"
    6 = "16.08%"
    7 = "-1.39%"
    8 = "2.07%"
    9 = "1.28%"
    10 = "3.26%"
    11 = "-1.64%"
    12 = "8.96%"
    13 = "-0.21%"
    14 = "0.78%"
    15 = "0.64%"
    16 = "-1.89%"
    17 = "1.61%"
    18 = "3.00%"
    19 = "0.16%"
    20 = "1.51%"
    21 = "-10.29%"
    22 = "7.34%"
    23 = "-0.49%"
    24 = "1.24%"
    25 = "0.98%"
    26 = "-7.49%"
    27 = "-0.44%"
    39 = "0.34%"
    40 = "1.81%"
    41 = "-5.60%"
    42 = "1.96%"
    43 = "-1.58%"
    44 = "-3.05%"
    45 = "-2.34%"
    46 = "0.02%"
    47 = "-15.98%"
    48 = "-8.86%"
    49 = "0.02%"
    50 = "0.09%"
}

foreach ($row in $dData.Keys) {
    $cell = $ws.Cells.Item($row, 4)   # column D
    $style = $cell.Style
    $cell.Value = "'" + $dData[$row]
    $cell.Style = $style
}

foreach ($row in $eData.Keys) {
    $cell = $ws.Cells.Item($row, 5)   # column E
    $style = $cell.Style
    $cell.Value = "'" + $eData[$row]
    $cell.Style = $style
}
